$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.805.49'
$ws.Range("E2").Value = '  -2.52%  '
$ws.Range("D3").Value = '2.447.19'
$ws.Range("E3").Value = '  -3.80%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.94%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.563'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0974'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.321'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.46%  '
$ws.Range("D13").Value = '2.883.71'
$ws.Range("E13").Value = '  -3.82%  '
$ws.Range("D14").Value = '57.741.64'
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.20%  '
$ws.Range("E16").Value = '  -3.35%  '
$ws.Range("D17").Value = '2.454.21'
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '312.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.401'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '2.567.26'
$ws.Range("E26").Value = '  -3.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.156'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.55%  '
$ws.Range("D30").Value = '0.0₃0733'
$ws.Range("E30").Value = '  -3.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.55%  '
$ws.Range("E33").Value = '  -8.27%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.94%  '
$ws.Range("E37").Value = '  -7.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.801'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.97%  '
$ws.Range("E41").Value = '  -5.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.581'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.77'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '255.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '122.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0918'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0489'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("E49").Value = '  -3.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.24%  '
